$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = 80
$ws.Range("E7").Value = 80
$ws.Range("F7").Value = 80
$ws.Range("I7").Value = 77.5
$ws.Range("K7").Value = 87.5
$ws.Range("L7").Value = 90
